# Scheduled-runner price/profit refresh: pushes newly-fetched market-board
# averages (and their downstream Leve price/profit columns) into the
# per-job Sheets. Cells are plain numeric literals (no formulas in this
# workbook), so each updated value is written directly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 156.81818
$ws.Range("I11").Value = 156.81818
$ws.Range("K11").Value = 156.81818
$ws.Range("M11").Value = -16.81818000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 561.2381
$ws.Range("I12").Value = 378
$ws.Range("J12").Value = 805.55554
$ws.Range("K12").Value = 378
$ws.Range("L12").Value = 805.55554
$ws.Range("M12").Value = -208
$ws.Range("N12").Value = -1145.55554

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2713.5715
$ws.Range("I62").Value = 2832.5
$ws.Range("K62").Value = 2832.5
$ws.Range("M62").Value = -2208.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2713.5715
$ws.Range("I65").Value = 2832.5
$ws.Range("K65").Value = 14162.5
$ws.Range("M65").Value = -11042.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3099.8
$ws.Range("I74").Value = 2666.3333
$ws.Range("J74").Value = 3750
$ws.Range("K74").Value = 2666.3333
$ws.Range("L74").Value = 3750
$ws.Range("M74").Value = -1730.3333
$ws.Range("N74").Value = -5622

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3099.8
$ws.Range("I77").Value = 2666.3333
$ws.Range("J77").Value = 3750
$ws.Range("K77").Value = 13331.6665
$ws.Range("L77").Value = 18750
$ws.Range("M77").Value = -8651.666499999999
$ws.Range("N77").Value = -28110

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1920
$ws.Range("I86").Value = 1933.6666
$ws.Range("K86").Value = 1933.6666
$ws.Range("M86").Value = -810.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1920
$ws.Range("I89").Value = 1933.6666
$ws.Range("K89").Value = 9668.333000000001
$ws.Range("M89").Value = -4052.333000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 2190
$ws.Range("J97").Value = 3900
$ws.Range("L97").Value = 11700
$ws.Range("N97").Value = -12692

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 4560.4707
$ws.Range("I106").Value = 3736.9092
$ws.Range("K106").Value = 3736.9092
$ws.Range("M106").Value = -3105.9092

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 61661.152
$ws.Range("J133").Value = 61661.152
$ws.Range("L133").Value = 61661.152
$ws.Range("N133").Value = -71781.152

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1932.6
$ws.Range("J137").Value = 2369.4285
$ws.Range("L137").Value = 7108.2855
$ws.Range("N137").Value = -12208.2855

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2160.0203
$ws.Range("I138").Value = 1438.0526
$ws.Range("J138").Value = 2331.4875
$ws.Range("K138").Value = 4314.1578
$ws.Range("L138").Value = 6994.462500000001
$ws.Range("M138").Value = 825.8422
$ws.Range("N138").Value = -17274.4625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3852.4592
$ws.Range("I32").Value = 3768.4639
$ws.Range("J32").Value = 12000
$ws.Range("K32").Value = 3768.4639
$ws.Range("L32").Value = 12000
$ws.Range("M32").Value = -3481.4639
$ws.Range("N32").Value = -12574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 44575.473
$ws.Range("I61").Value = 67157.164
$ws.Range("J61").Value = 5864
$ws.Range("K61").Value = 67157.164
$ws.Range("L61").Value = 5864
$ws.Range("M61").Value = -66945.164
$ws.Range("N61").Value = -6288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 757.4820999999999
$ws.Range("I74").Value = 651.1731
$ws.Range("K74").Value = 651.1731
$ws.Range("M74").Value = 222.8269

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 757.4820999999999
$ws.Range("I77").Value = 651.1731
$ws.Range("K77").Value = 3255.8655
$ws.Range("M77").Value = 1112.1345

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 44575.473
$ws.Range("I136").Value = 67157.164
$ws.Range("J136").Value = 5864
$ws.Range("K136").Value = 201471.492
$ws.Range("L136").Value = 17592
$ws.Range("M136").Value = -198921.492
$ws.Range("N136").Value = -22692

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 292
$ws.Range("I94").Value = 322.66666
$ws.Range("K94").Value = 322.66666
$ws.Range("M94").Value = 128.33334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2006.3334
$ws.Range("I99").Value = 2007.5
$ws.Range("J99").Value = 2004
$ws.Range("K99").Value = 2007.5
$ws.Range("L99").Value = 2004
$ws.Range("M99").Value = -509.5
$ws.Range("N99").Value = -5000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2142.3845
$ws.Range("I105").Value = 2095.8286
$ws.Range("J105").Value = 2549.75
$ws.Range("K105").Value = 2095.8286
$ws.Range("L105").Value = 2549.75
$ws.Range("M105").Value = -348.8285999999998
$ws.Range("N105").Value = -6043.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3927
$ws.Range("I107").Value = 3927
$ws.Range("K107").Value = 3927
$ws.Range("M107").Value = -2007

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5143.7407
$ws.Range("I134").Value = 4905.273
$ws.Range("J134").Value = 6193
$ws.Range("K134").Value = 14715.819
$ws.Range("L134").Value = 18579
$ws.Range("M134").Value = -12180.819
$ws.Range("N134").Value = -23649

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1319390.6
$ws.Range("J58").Value = 1740.24
$ws.Range("L58").Value = 1740.24
$ws.Range("N58").Value = -2146.24

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1562.7084
$ws.Range("I132").Value = 1060.4
$ws.Range("K132").Value = 3181.2
$ws.Range("M132").Value = -651.2000000000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1418.2222
$ws.Range("I134").Value = 1276.9546
$ws.Range("K134").Value = 3830.8638
$ws.Range("M134").Value = -1295.8638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1319390.6
$ws.Range("J136").Value = 1740.24
$ws.Range("L136").Value = 5220.72
$ws.Range("N136").Value = -10320.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 389.81818
$ws.Range("J38").Value = 490.14285
$ws.Range("L38").Value = 1470.42855
$ws.Range("N38").Value = -2164.42855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1207445.2
$ws.Range("I132").Value = 1544310.4
$ws.Range("K132").Value = 4632931.199999999
$ws.Range("M132").Value = -4630401.199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2470.5
$ws.Range("I7").Value = 1985.6666
$ws.Range("J7").Value = 3925
$ws.Range("K7").Value = 1985.6666
$ws.Range("L7").Value = 3925
$ws.Range("M7").Value = -1873.6666
$ws.Range("N7").Value = -4149

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2227.6667
$ws.Range("I93").Value = 1063.8334
$ws.Range("K93").Value = 1063.8334
$ws.Range("M93").Value = 184.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2108.4
$ws.Range("I100").Value = 1130.8334
$ws.Range("K100").Value = 1130.8334
$ws.Range("M100").Value = -589.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2470.5
$ws.Range("I126").Value = 1985.6666
$ws.Range("J126").Value = 3925
$ws.Range("K126").Value = 5956.9998
$ws.Range("L126").Value = 11775
$ws.Range("M126").Value = -3486.9998
$ws.Range("N126").Value = -16715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3675.4285
$ws.Range("I132").Value = 1923.3
$ws.Range("K132").Value = 5769.9
$ws.Range("M132").Value = -3239.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2862.75
$ws.Range("I136").Value = 2846.2144
$ws.Range("K136").Value = 8538.643199999999
$ws.Range("M136").Value = -5988.643199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 12990.286
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 137996.17
$ws.Range("I122").Value = 163595.4
$ws.Range("K122").Value = 490786.2
$ws.Range("M122").Value = -488336.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6510.773
$ws.Range("I126").Value = 7686.375
$ws.Range("K126").Value = 23059.125
$ws.Range("M126").Value = -20589.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1632.0358
$ws.Range("I132").Value = 1468.3096
$ws.Range("J132").Value = 2123.2144
$ws.Range("K132").Value = 4404.9288
$ws.Range("L132").Value = 6369.6432
$ws.Range("M132").Value = -1874.9288
$ws.Range("N132").Value = -11429.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 15874546
$ws.Range("I136").Value = 23149214
$ws.Range("K136").Value = 69447642
$ws.Range("M136").Value = -69445092
